# Error Calculations and Plots
# Apply the missing-data imputation / correction changes to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C single-cell corrections (values swapped between present/missing) ---
$ws.Range("C6").Value = 15.1
$ws.Range("C8").Value = ""
$ws.Range("C12").Value = 12.5
$ws.Range("C14").Value = ""
$ws.Range("C17").Value = 11.2
$ws.Range("C18").Value = 11.5
$ws.Range("C19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("C23").Value = 12.2

# --- Remove the "RM 232" row (row 26) entirely ---
$ws.Rows(26).Delete()

# --- Remove the "SC 92" row (now row 27 after the previous delete) entirely ---
$ws.Rows(27).Delete()

# --- Column B/C corrections on the rows that shifted up ---
# SC 101 (row 27): B was blank -> -20.4 ; C was 10 -> blank
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""
# SC 105 (row 28): B was -19.6 -> blank
$ws.Range("B28").Value = ""
# SC 119 (row 29): B was -19.5 -> blank
$ws.Range("B29").Value = ""
# SC 120 (row 30): B was blank -> -19.7
$ws.Range("B30").Value = -19.7
# SC 193 (row 32): B was -19.9 -> blank
$ws.Range("B32").Value = ""
